$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.071.42'
$ws.Range("E2").Value = '  -3.00%  '

$ws.Range("D3").Value = '1.653.04'
$ws.Range("E3").Value = '  -4.78%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.91'
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("E7").Value = '  -7.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '39.53'
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2622'
$ws.Range("E9").Value = '  -4.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.05988'
$ws.Range("E10").Value = '  -2.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.07091'
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.661.17'
$ws.Range("E12").Value = '  -4.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '14.43'
$ws.Range("E13").Value = '  -3.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.6188'
$ws.Range("E14").Value = '  -3.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '4.581'
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '72.98'
$ws.Range("E16").Value = '  -5.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D17").Value = '0.9999'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '25.067.09'
$ws.Range("E19").Value = '  -3.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '11.37'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.000006572'
$ws.Range("E21").Value = '  -2.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.428'
$ws.Range("E22").Value = '  +4.01%  '

$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '1.867.07'
$ws.Range("E23").Value = '  -4.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '8.455'
$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '5.246'
$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '132.65'
$ws.Range("E26").Value = '  -3.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '14.71'
$ws.Range("E27").Value = '  -3.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '1.392'
$ws.Range("E28").Value = '  -7.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '1.695'
$ws.Range("E29").Value = '  -3.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '101.42'
$ws.Range("E30").Value = '  -3.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '3.795'
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.07911'
$ws.Range("E32").Value = '  -3.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '3.523'
$ws.Range("E33").Value = '  -3.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.04587'
$ws.Range("E34").Value = '  -0.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.606'
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9381'
$ws.Range("E36").Value = '  -4.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.5781'
$ws.Range("E37").Value = '  -6.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.625'
$ws.Range("E38").Value = '  -2.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01535'
$ws.Range("E39").Value = '  -3.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.8418'
$ws.Range("E40").Value = '  +13.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '1.826'
$ws.Range("E42").Value = '  -4.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '98.68'
$ws.Range("E43").Value = '  -1.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3703'
$ws.Range("E44").Value = '  -3.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '4.827'
$ws.Range("E45").Value = '  -3.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1114'
$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '6.055'
$ws.Range("E47").Value = '  -2.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05146'
$ws.Range("E48").Value = '  -1.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '52.07'
$ws.Range("E49").Value = '  -4.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '29.35'
$ws.Range("E50").Value = '  -3.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  -0.07%  '
